# BOM.xlsx edit: remove the RSMF2JT100R resistor line, add CF14JT100R and
# CF18JT3K00 resistor lines, and keep the rest of the BOM (KEMET, Nichicon,
# totals, footnotes) intact below them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the old "Stackpole Electronics RSMF2JT100R" row (row 13).
#    This shifts CF12JT10R0 / KEMET / Nichicon / totals / footnotes up by one.
$ws.Rows.Item(13).Delete()

# 2) Make room for the two new resistor lines right after CF12JT10R0
#    (which is now row 13), pushing KEMET / Nichicon / totals / footnotes
#    back down by two rows.
$ws.Rows.Item(14).Resize(2).Insert()

# 3) Update the CF12JT10R0 line (row 13) quantities.
$ws.Cells.Item(13,3).Value = 0.1
$ws.Cells.Item(13,4).Value = 2
$ws.Cells.Item(13,5).Formula = "=CEILING.MATH(3*D13)"
$ws.Cells.Item(13,6).Formula = "=C13*D13"

# 4) Fill in the new CF14JT100R line (row 14).
$ws.Cells.Item(14,1).Value = "Stackpole Electronics  CF14JT100R "
$ws.Cells.Item(14,2).Value = "Through Hole 100 ohm ¼ watt Resistor"
$ws.Cells.Item(14,3).Value = 0.1
$ws.Cells.Item(14,4).Value = 14
$ws.Cells.Item(14,5).Formula = "=CEILING.MATH(3*D14)"
$ws.Cells.Item(14,6).Formula = "=C14*D14"
$ws.Cells.Item(14,7).Value = "Digi-key"
$ws.Cells.Item(14,8).Value = "https://www.digikey.com/product-detail/en/stackpole-electronics-inc/CF14JT100R/CF14JT100RCT-ND/1830327"

# 5) Fill in the new CF18JT3K00 line (row 15).
$ws.Cells.Item(15,1).Value = "Stackpole Electronics  CF18JT3K00 "
$ws.Cells.Item(15,2).Value = "Through Hole 3k ohm 1/4 watt Resistor"
$ws.Cells.Item(15,3).Value = 0.1
$ws.Cells.Item(15,4).Value = 14
$ws.Cells.Item(15,5).Formula = "=CEILING.MATH(3*D15)"
$ws.Cells.Item(15,6).Formula = "=C15*D15"
$ws.Cells.Item(15,7).Value = "Digi-key"
$ws.Cells.Item(15,8).Value = "https://www.digikey.com/product-detail/en/stackpole-electronics-inc/CF18JT3K00/CF18JT3K00CT-ND/2022753"

# 6) KEMET (row 16) and Nichicon (row 17) only moved down two rows -- their
#    data followed automatically, but re-assert their E/F formulas so they
#    stay plain formulas (rather than picking up an implicit-intersection
#    wrapper from the row shift).
$ws.Cells.Item(16,5).Formula = "=CEILING.MATH(3*D16)"
$ws.Cells.Item(16,6).Formula = "=C16*D16"
$ws.Cells.Item(17,5).Formula = "=CEILING.MATH(3*D17)"
$ws.Cells.Item(17,6).Formula = "=C17*D17"

# 7) The "Total Cost per Board" sum needs to cover the newly-added F17
#    (Nichicon, after the shift) as well.
$ws.Cells.Item(21,6).Formula = "=F2+F3+F4+F5+F6+F7+F8+F9+F10+F11+F12+F13+F14+F15+F16+F17"

# 8) Keep the selection / active cell in sync with the new layout.
$ws.Range("B23").Select()
